# Daily attendance processing - 2026-01-03 20:58:43
# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever both are present, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
